$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New QualityProfile column block, entered in the same order the
#     author actually typed it (this governs the shared-string table order) ---
$ws.Range("E2").Value = "QualityProfile"

$ws.Range("E14").Value = "Quality is pretty shitty…"
$ws.Range("F14").Value = 200
$ws.Range("G14").Value = 250

$ws.Range("F2").Value = "F-truncation"
$ws.Range("G2").Value = "R-truncation"

$ws.Range("E12").Value = "R2 is pretty shitty"
$ws.Range("F12").Value = 150
$ws.Range("G12").Value = 200

$ws.Range("E13").Value = "Not too bad"
$ws.Range("F13").Value = 125
$ws.Range("G13").Value = 100

$ws.Range("E11").Value = "R1 is unusually terrible"
$ws.Range("F11").Value = 150
$ws.Range("G11").Value = 100

$ws.Range("E8").Value = "R1 is good, R2 less so"
$ws.Range("F8").Value = 225
$ws.Range("G8").Value = 250

$ws.Range("E5").Value = "OK"
$ws.Range("F5").Value = 150
$ws.Range("G5").Value = 200

$ws.Range("E6").Value = "OK"
$ws.Range("F6").Value = 200
$ws.Range("G6").Value = 225

# --- Mark TMM_DSS (rows 7-8) and MS_DSS (rows 9-10) as confirmed,
#     matching the highlight formatting already used for rows 3-4 ---
$ws.Range("A4:D4").Copy() | Out-Null
$ws.Range("A7:D7").PasteSpecial(-4122) | Out-Null
$ws.Range("A9:D10").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

# --- Column widths for the new columns ---
$ws.Columns("E").ColumnWidth = 20.17
$ws.Columns("G").ColumnWidth = 11

# --- Update the saved cursor/selection position ---
$ws.Range("F24").Select() | Out-Null
